$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.728861689567566
$ws.Range("B1").Value = 1.904061675071716
$ws.Range("C1").Value = 2.160602331161499
$ws.Range("D1").Value = 2.681890487670898
$ws.Range("E1").Value = 1.532042026519775
